# Tasks.xlsx - Added new team tasks
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-DoneStyle($cell) {
    $cell.Interior.Color = 5296274   # RGB(146,208,80) green
    $cell.WrapText = $true
}

function Set-PlainStyle($cell) {
    # plain wrap-only style (style index 1), no fill
    $cell.WrapText = $true
}

# ---- Row 15: F15 -> Done (green) ----
$f15 = $ws.Cells.Item(15, 6)
$f15.Value = "Done"
Set-DoneStyle $f15

# ---- Row 22: F22 -> Done (green); G22 -> 2 hour ----
$f22 = $ws.Cells.Item(22, 6)
$f22.Value = "Done"
Set-DoneStyle $f22
$ws.Cells.Item(22, 7).Value = "2 hour"

# ---- Row 23: clear out old task-17 row content, leave only "No" (A23=20) ----
$ws.Range("F23:G23").Clear()
$ws.Range("B23:E23").ClearContents()

# ---- Row 24: Student-Support Solutions app development ----
$ws.Cells.Item(24, 4).Value = "Tasks 1-20 are finished"
$ws.Cells.Item(24, 2).Value = "Student-Support Solutions app development"
$ws.Range("C24").ClearContents()

# ---- Row 25: Student-My tickets app development ----
$ws.Cells.Item(25, 2).Value = "Student-My tickets app development"
$ws.Range("C25").ClearContents()
$ws.Cells.Item(25, 4).Value = "Tasks 1-20 are finished"

# ---- Row 26: Student-New Ticket app development ----
$ws.Cells.Item(26, 2).Value = "Student-New Ticket app development"
$ws.Range("C26").ClearContents()
$ws.Cells.Item(26, 4).Value = "Tasks 1-20 are finished"

# ---- Row 27: Support team Member-Solution Categoryapp development ----
$ws.Cells.Item(27, 2).Value = "Support team Member-Solution Categoryapp development"
$ws.Cells.Item(27, 3).Value = "Shamil"
$ws.Cells.Item(27, 4).Value = "Tasks 1-20 are finished"

# ---- New Row 28 ----
$a28 = $ws.Cells.Item(28, 1)
$a28.Value = 25
Set-PlainStyle $a28
$b28 = $ws.Cells.Item(28, 2)
$b28.Value = "Support team Member-Support Solutions app development"
Set-PlainStyle $b28
$d28 = $ws.Cells.Item(28, 4)
$d28.Value = "Tasks 1-20 are finished"
Set-PlainStyle $d28
$ws.Cells.Item(28, 6).Value = "New"

# ---- New Row 29 ----
$a29 = $ws.Cells.Item(29, 1)
$a29.Value = 26
Set-PlainStyle $a29
$b29 = $ws.Cells.Item(29, 2)
$b29.Value = "Support team Member-Open tickets app development"
Set-PlainStyle $b29
$d29 = $ws.Cells.Item(29, 4)
$d29.Value = "Tasks 1-20 are finished"
Set-PlainStyle $d29
$ws.Cells.Item(29, 6).Value = "New"

# ---- New Row 30 ----
$a30 = $ws.Cells.Item(30, 1)
$a30.Value = 27
Set-PlainStyle $a30
$b30 = $ws.Cells.Item(30, 2)
$b30.Value = "Support team Manager-Assigne Responsible app development"
Set-PlainStyle $b30
$d30 = $ws.Cells.Item(30, 4)
$d30.Value = "Tasks 1-20 are finished"
Set-PlainStyle $d30
$ws.Cells.Item(30, 6).Value = "New"

# ---- New Row 31 ----
$a31 = $ws.Cells.Item(31, 1)
$a31.Value = 28
Set-PlainStyle $a31
$b31 = $ws.Cells.Item(31, 2)
$b31.Value = "Support team Manager-Analytics app development"
Set-PlainStyle $b31
$c31 = $ws.Cells.Item(31, 3)
$c31.Value = "Shamil"
Set-PlainStyle $c31
$d31 = $ws.Cells.Item(31, 4)
$d31.Value = "Tasks 1-20 are finished"
Set-PlainStyle $d31
$ws.Cells.Item(31, 6).Value = "New"

# ---- Update view: selection + scroll position ----
$ws.Range("B31").Select()
$excel.ActiveWindow.ScrollRow = 13
$excel.ActiveWindow.ScrollColumn = 1
